$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume/issue number, week-of date range) ---
$ws.Range('A8').Value = 'Volume 32   Number  10'
$ws.Range('C9').Value = 'Report Covering the Week  3/3/2025  Through  3/9/2025'

# --- Cells that become text placeholders ("0") - copy format+value from C14 which already holds it ---
$ws.Range('C14').Copy($ws.Range('C15'))
$ws.Range('C14').Copy($ws.Range('D15'))
$ws.Range('C14').Copy($ws.Range('D23'))
$ws.Range('C14').Copy($ws.Range('C27'))
$ws.Range('C14').Copy($ws.Range('D29'))
$ws.Range('C14').Copy($ws.Range('D30'))

# --- Cells that become text placeholders ("***.*") - copy format+value from E14 which already holds it ---
$ws.Range('E14').Copy($ws.Range('E15'))
$ws.Range('E14').Copy($ws.Range('E23'))
$ws.Range('E14').Copy($ws.Range('E29'))
$ws.Range('E14').Copy($ws.Range('E30'))

# --- Numeric value updates ---
$ws.Range('F15').Value = 4
$ws.Range('G15').Value = 3
$ws.Range('H15').Value = 33.333333333333
$ws.Range('N15').Value = 75
$ws.Range('C16').Value = 2
$ws.Range('D16').Value = 6
$ws.Range('E16').Value = -66.666666666666
$ws.Range('F16').Value = 11
$ws.Range('G16').Value = 19
$ws.Range('H16').Value = -42.105263157894
$ws.Range('I16').Value = 31
$ws.Range('J16').Value = 51
$ws.Range('K16').Value = -39.215686274509
$ws.Range('L16').Value = -35.416666666666
$ws.Range('M16').Value = -20.51282051282
$ws.Range('N16').Value = -77.697841726618
$ws.Range('C17').Value = 5
$ws.Range('D17').Value = 9
$ws.Range('E17').Value = -44.444444444444
$ws.Range('F17').Value = 31
$ws.Range('G17').Value = 34
$ws.Range('H17').Value = -8.823529411764
$ws.Range('I17').Value = 79
$ws.Range('J17').Value = 84
$ws.Range('K17').Value = -5.95238095238
$ws.Range('L17').Value = 29.508196721311
$ws.Range('M17').Value = 3.947368421052
$ws.Range('N17').Value = 19.696969696969
$ws.Range('C18').Value = 6
$ws.Range('D18').Value = 2
$ws.Range('E18').Value = 200
$ws.Range('F18').Value = 21
$ws.Range('G18').Value = 11
$ws.Range('H18').Value = 90.90909090909
$ws.Range('I18').Value = 39
$ws.Range('J18').Value = 29
$ws.Range('K18').Value = 34.482758620689
$ws.Range('L18').Value = 5.405405405405
$ws.Range('M18').Value = 5.405405405405
$ws.Range('N18').Value = -74
$ws.Range('C19').Value = 3
$ws.Range('D19').Value = 6
$ws.Range('E19').Value = -50
$ws.Range('F19').Value = 13
$ws.Range('G19').Value = 31
$ws.Range('H19').Value = -58.064516129032
$ws.Range('I19').Value = 54
$ws.Range('J19').Value = 72
$ws.Range('K19').Value = -25
$ws.Range('L19').Value = -34.939759036144
$ws.Range('M19').Value = 14.893617021276
$ws.Range('N19').Value = -14.285714285714
$ws.Range('C20').Value = 2
$ws.Range('D20').Value = 2
$ws.Range('E20').Value = 0
$ws.Range('I20').Value = 28
$ws.Range('J20').Value = 31
$ws.Range('K20').Value = -9.677419354838
$ws.Range('L20').Value = -46.153846153846
$ws.Range('M20').Value = 86.666666666666
$ws.Range('N20').Value = -59.420289855072
$ws.Range('C21').Value = 18
$ws.Range('D21').Value = 25
$ws.Range('E21').Value = -28
$ws.Range('F21').Value = 91
$ws.Range('G21').Value = 108
$ws.Range('H21').Value = -15.74074074074
$ws.Range('I21').Value = 238
$ws.Range('J21').Value = 275
$ws.Range('K21').Value = -13.454545454545
$ws.Range('L21').Value = -17.931034482758
$ws.Range('M21').Value = 10.697674418604
$ws.Range('N21').Value = -52.208835341365
$ws.Range('G22').Value = 4
$ws.Range('J22').Value = 6
$ws.Range('K22').Value = -50
$ws.Range('C24').Value = 24
$ws.Range('D24').Value = 18
$ws.Range('E24').Value = 33.333333333333
$ws.Range('F24').Value = 81
$ws.Range('G24').Value = 63
$ws.Range('H24').Value = 28.571428571428
$ws.Range('I24').Value = 166
$ws.Range('J24').Value = 139
$ws.Range('K24').Value = 19.424460431654
$ws.Range('L24').Value = 24.812030075188
$ws.Range('M24').Value = 90.804597701149
$ws.Range('C25').Value = 11
$ws.Range('D25').Value = 4
$ws.Range('E25').Value = 175
$ws.Range('F25').Value = 25
$ws.Range('G25').Value = 20
$ws.Range('H25').Value = 25
$ws.Range('I25').Value = 52
$ws.Range('J25').Value = 42
$ws.Range('K25').Value = 23.809523809523
$ws.Range('L25').Value = 26.829268292682
$ws.Range('C26').Value = 11
$ws.Range('D26').Value = 9
$ws.Range('E26').Value = 22.222222222222
$ws.Range('F26').Value = 50
$ws.Range('G26').Value = 44
$ws.Range('H26').Value = 13.636363636363
$ws.Range('I26').Value = 111
$ws.Range('J26').Value = 117
$ws.Range('K26').Value = -5.128205128205
$ws.Range('L26').Value = 38.75
$ws.Range('M26').Value = 13.265306122449
$ws.Range('E27').Value = -100
$ws.Range('F27').Value = 5
$ws.Range('H27').Value = 0
$ws.Range('J27').Value = 13
$ws.Range('K27').Value = -15.384615384615
$ws.Range('C28').Value = 2
$ws.Range('E28').Value = 0
$ws.Range('F28').Value = 11
$ws.Range('G28').Value = 8
$ws.Range('H28').Value = 37.5
$ws.Range('I28').Value = 21
$ws.Range('J28').Value = 18
$ws.Range('K28').Value = 16.666666666666
$ws.Range('L28').Value = 0
